# Update countries & provincias Spain
# Applies the COVID dashboard data refresh captured by the diff:
#   - Updated "Datos actualizados" timestamp (08:07 -> 09:24)
#   - Refreshed numeric figures (Casos totales / Nuevos casos / Casos activos /
#     Recuperados / Casos criticos / Muertes hoy / Muertes) for a handful of
#     countries (Israel, Nigeria, Austria, Armenia, Hungria, Montenegro,
#     Georgia, Gibraltar, Taiwan).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp footer (row 1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 09:24"

# --- Numeric refresh (columns B:H, row headers: Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) -------

function Set-Row($r, $vals) {
    $col = 2
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col++
    }
}

# Row 27: Israel
Set-Row 27 @(301289, 1088, 260803, 38358, 0, 1, 2128)

# Row 61: Nigeria
Set-Row 61 @(61460, 1465, 47541, 12863, 0, 10, 1056)

# Row 62: Austria
Set-Row 62 @(60982, 0, 52194, 7672, 0, 0, 1116)

# Row 63: Armenia
Set-Row 63 @(60224, 0, 46798, 12549, 0, 0, 877)

# Row 75: Hungria
Set-Row 75 @(43025, 1293, 13134, 28806, 0, 33, 1085)

# Row 99: Montenegro
Set-Row 99 @(15327, 887, 7613, 7590, 0, 11, 124)

# Row 100: Georgia
Set-Row 100 @(14672, 0, 10355, 4096, 0, 0, 221)

# Row 176: Gibraltar
Set-Row 176 @(535, 4, 491, 37, 0, 0, 7)

# Row 177: Taiwan
Set-Row 177 @(531, 0, 439, 92, 0, 0, 0)
